$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the first
#    paragraph (the Heading1 title "Play Big Catch Slot for Free - Game
#    Review").
# ----------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titleEnd = $titlePara.Range
$titleEnd.Collapse(0)
$titleEnd.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaRng = $metaPara.Range

$metaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body><w:p>' +
  '<w:r/>' +
  '<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>' +
  '<w:r><w:t>: Read our review of Big Catch, a simple yet engaging slot game by Novomatic with free spins and bonus features. Play for free and catch the big one!</w:t></w:r>' +
  '</w:p></w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

[void]$metaRng.InsertXML($metaXml)

# ----------------------------------------------------------------------
# 2) Near the end of the document, remove the duplicated bold heading
#    paragraph ("Play Big Catch Slot for Free - Game Review") and turn
#    the following italic paragraph's meta-description text into the
#    image-generation prompt text.
# ----------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 2; $i--) {
  $p = $d.Paragraphs($i)
  $txt = $p.Range.Text.TrimEnd([char]13, [char]7)
  if ($txt -eq "Play Big Catch Slot for Free - Game Review") {
    [void]$p.Range.Delete()
    break
  }
}

$count = $d.Paragraphs.Count
$promptPara = $d.Paragraphs($count)
$promptRng = $promptPara.Range
$target = $d.Range($promptRng.Start, $promptRng.End - 1)

$promptText = 'Prompt: Create an eye-catching feature image for the game "Big Catch". The image should be in cartoon style and feature a happy Maya warrior with glasses. The warrior can be holding a big hook and surrounded by sea creatures. The overall theme of the image should be underwater adventure and fishing. Use bold and vibrant colors to draw in the viewer''s attention and make the image stand out.'

$promptXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body><w:p>' +
  '<w:r><w:rPr><w:i/></w:rPr><w:t>' + $promptText + '</w:t></w:r>' +
  '</w:p></w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

[void]$target.InsertXML($promptXml)
